$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Entities")

# Fix Spawn Room / Coin value: D3 300 -> 10
$ws.Range("D3").Value = 10

# Fix Weapon stats
# DarkSword (row 5): damage 10 -> 20, rate 0.8 -> 1.3
$ws.Range("G5").Value = 20
$ws.Range("H5").Value = 1.3

# Sword (row 6): damage 20 -> 10, rate 0.4 -> 0.8
$ws.Range("G6").Value = 10
$ws.Range("H6").Value = 0.8

# ChainSaw (row 7): damage 30 -> 50, rate 1.5 -> 2
$ws.Range("G7").Value = 50
$ws.Range("H7").Value = 2

# Update selected cell to match author's last selection
$ws.Range("G6").Select()
